$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "EC3268"
$ws.Range("B4").Value = "Sachin Lonkar"
$ws.Range("C4").Value = "B"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = ""
